# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# table to the refreshed figures from the latest scrape.
#
# Numeric-looking Price values (e.g. "95.82") are written with a
# temporary "@" (text) number format so Excel keeps them as literal
# strings instead of silently parsing them into floating point numbers
# - the workbook always stores these as plain text. The style is reset
# back to "Normal" immediately afterwards so no stray cell formatting
# is left behind. Price values that already contain multiple separators
# (e.g. "42.222.98") or subscript digits are unambiguous text and are
# assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.222.98"
$ws.Range("E2").Value = "  -1.48%  "
$ws.Range("D3").Value = "2.272.32"
$ws.Range("E3").Value = "  -2.36%  "
$ws.Range("E4").Value = "  +0.04%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "299.06"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -2.29%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "95.82"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -4.72%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.495"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -2.60%  "
$ws.Range("E8").Value = "  +0.06%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.495"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -1.99%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "33.45"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -2.94%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0792"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -0.22%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "48.17"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -7.70%  "
$ws.Range("E13").Value = "  -0.12%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "6.68"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.05%  "
$ws.Range("D15").Value = "2.623.79"
$ws.Range("E15").Value = "  -2.70%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "15.52"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("D17").Value = "2.277.68"
$ws.Range("E17").Value = "  -2.57%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.783"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -4.79%  "
$ws.Range("D19").Value = "42.155.34"
$ws.Range("E19").Value = "  -1.46%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "11.76"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +1.64%  "
$ws.Range("D21").Value = "0.0₃0893"
$ws.Range("E21").Value = "  -1.38%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "6.00"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -2.41%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "66.70"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -3.38%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "235.32"
$cell.Style = "Normal"
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "1.97"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("E26").Value = "  +0.29%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "2.46"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -3.36%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "23.99"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -5.33%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.18"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -2.32%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "168.59"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +4.69%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "34.10"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -1.78%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "9.14"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -0.72%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -0.20%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "4.92"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -2.76%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "4.53"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -1.53%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "16.60"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -2.99%  "
$ws.Range("E37").Value = "  -4.85%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.0688"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -4.23%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.79"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -3.23%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.0988"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -2.34%  "
$ws.Range("E41").Value = "  -2.47%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "1.73"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -5.29%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "2.46"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -4.46%  "
$ws.Range("D44").Value = "1.960.97"
$ws.Range("E44").Value = "  -2.44%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.0278"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -1.63%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "17.55"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -5.07%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "9.55"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -5.86%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "2.80"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -2.85%  "
$ws.Range("D49").Value = "2.495.98"
$ws.Range("E49").Value = "  -2.34%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "52.38"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -5.54%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "4.54"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -3.09%  "
